$d = $word.ActiveDocument

# The edit removes the "+" in every "B+" occurrence (turning "B+ tree"/"B+ trees"
# into "B tree"/"B trees"). In the saved OOXML this also leaves the lone "B"
# isolated into its own run (three runs with identical formatting: the text
# before, the single "B", and the text after), matching how Word naturally
# fragments a run when a sub-range inside it receives a direct formatting
# touch.
#
# Phase 1: delete every stray "+" first. (Deleting text causes the host to
# re-coalesce same-format runs within the edited paragraph, so any run split
# performed before all the deletions are done would get merged away again.)
$bStarts = New-Object System.Collections.ArrayList
$cont = $true
while ($cont) {
    $rng = $d.Content
    $found = $rng.Find.Execute("B+", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $start = $rng.Start
        [void]$bStarts.Add($start)
        $plusRng = $d.Range($start + 1, $start + 2)
        $plusRng.Delete()
    } else {
        $cont = $false
    }
}

# Phase 2: now that all the text edits are done, isolate each remaining lone
# "B" into its own run via a no-op formatting round-trip.
foreach ($start in $bStarts) {
    $bRng = $d.Range($start, $start + 1)
    $bRng.Bold = 1
    $bRng.Bold = 0
}
